# Auto-generated edit script applying numeric updates to the Famfrit_Profits workbook
# (8 sheets: ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(5, 8).Value = 907.7143  # H5
$ws.Cells.Item(5, 9).Value = 317.5  # I5
$ws.Cells.Item(5, 10).Value = 1694.6666  # J5
$ws.Cells.Item(5, 11).Value = 317.5  # K5
$ws.Cells.Item(5, 12).Value = 1694.6666  # L5
$ws.Cells.Item(5, 13).Value = -202.5  # M5
$ws.Cells.Item(5, 14).Value = -1924.6666  # N5
$ws.Cells.Item(43, 8).Value = 1179.5454  # H43
$ws.Cells.Item(43, 9).Value = 1039.2858  # I43
$ws.Cells.Item(43, 11).Value = 1039.2858  # K43
$ws.Cells.Item(43, 13).Value = -970.2858000000001  # M43
$ws.Cells.Item(62, 8).Value = 5874.875  # H62
$ws.Cells.Item(62, 10).Value = 5874.875  # J62
$ws.Cells.Item(62, 12).Value = 5874.875  # L62
$ws.Cells.Item(62, 14).Value = -7122.875  # N62
$ws.Cells.Item(65, 8).Value = 5874.875  # H65
$ws.Cells.Item(65, 10).Value = 5874.875  # J65
$ws.Cells.Item(65, 12).Value = 29374.375  # L65
$ws.Cells.Item(65, 14).Value = -35614.375  # N65
$ws.Cells.Item(98, 8).Value = 4064.639  # H98
$ws.Cells.Item(98, 9).Value = 3354.7778  # I98
$ws.Cells.Item(98, 11).Value = 3354.7778  # K98
$ws.Cells.Item(98, 13).Value = -1856.7778  # M98
$ws.Cells.Item(122, 8).Value = 4064.639  # H122
$ws.Cells.Item(122, 9).Value = 3354.7778  # I122
$ws.Cells.Item(122, 11).Value = 10064.3334  # K122
$ws.Cells.Item(122, 13).Value = -7614.3334  # M122
$ws.Cells.Item(132, 8).Value = 3745.7778  # H132
$ws.Cells.Item(132, 9).Value = 3365.44  # I132
$ws.Cells.Item(132, 11).Value = 10096.32  # K132
$ws.Cells.Item(132, 13).Value = -7566.32  # M132
$ws.Cells.Item(137, 8).Value = 4230.069  # H137
$ws.Cells.Item(137, 9).Value = 1505.2916  # I137
$ws.Cells.Item(137, 10).Value = 17309  # J137
$ws.Cells.Item(137, 11).Value = 4515.8748  # K137
$ws.Cells.Item(137, 12).Value = 51927  # L137
$ws.Cells.Item(137, 13).Value = -1965.8748  # M137
$ws.Cells.Item(137, 14).Value = -57027  # N137

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 3854.8022  # H32
$ws.Cells.Item(32, 9).Value = 2997.0132  # I32
$ws.Cells.Item(32, 11).Value = 2997.0132  # K32
$ws.Cells.Item(32, 13).Value = -2710.0132  # M32
$ws.Cells.Item(122, 8).Value = 20836444  # H122
$ws.Cells.Item(122, 9).Value = 2341.1428  # I122
$ws.Cells.Item(122, 10).Value = 37040748  # J122
$ws.Cells.Item(122, 11).Value = 7023.428400000001  # K122
$ws.Cells.Item(122, 12).Value = 111122244  # L122
$ws.Cells.Item(122, 13).Value = -4573.428400000001  # M122
$ws.Cells.Item(122, 14).Value = -111127144  # N122

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(99, 8).Value = 2732.111  # H99
$ws.Cells.Item(99, 9).Value = 1811.7273  # I99
$ws.Cells.Item(99, 11).Value = 1811.7273  # K99
$ws.Cells.Item(99, 13).Value = -313.7273  # M99
$ws.Cells.Item(107, 8).Value = 2934.8462  # H107
$ws.Cells.Item(107, 9).Value = 2346.1667  # I107
$ws.Cells.Item(107, 11).Value = 2346.1667  # K107
$ws.Cells.Item(107, 13).Value = -426.1667000000002  # M107
$ws.Cells.Item(112, 8).Value = 0  # H112
$ws.Cells.Item(112, 9).Value = 0  # I112
$ws.Cells.Item(112, 10).Value = 0  # J112
$ws.Cells.Item(112, 11).Value = 0  # K112
$ws.Cells.Item(112, 12).Value = 0  # L112
$ws.Cells.Item(112, 13).ClearContents()  # M112
$ws.Cells.Item(112, 14).ClearContents()  # N112

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(7, 8).Value = 285  # H7
$ws.Cells.Item(7, 9).Value = 200  # I7
$ws.Cells.Item(7, 10).Value = 299.16666  # J7
$ws.Cells.Item(7, 11).Value = 200  # K7
$ws.Cells.Item(7, 12).Value = 299.16666  # L7
$ws.Cells.Item(7, 13).Value = -87  # M7
$ws.Cells.Item(7, 14).Value = -525.16666  # N7
$ws.Cells.Item(22, 8).Value = 6175.6313  # H22
$ws.Cells.Item(22, 9).Value = 8042.923  # I22
$ws.Cells.Item(22, 10).Value = 2129.8333  # J22
$ws.Cells.Item(22, 11).Value = 8042.923  # K22
$ws.Cells.Item(22, 12).Value = 2129.8333  # L22
$ws.Cells.Item(22, 13).Value = -7692.923  # M22
$ws.Cells.Item(22, 14).Value = -2829.8333  # N22
$ws.Cells.Item(57, 8).Value = 44495  # H57
$ws.Cells.Item(57, 10).Value = 44495  # J57
$ws.Cells.Item(57, 12).Value = 44495  # L57
$ws.Cells.Item(57, 14).Value = -45615  # N57
$ws.Cells.Item(86, 8).Value = 13858  # H86
$ws.Cells.Item(86, 10).Value = 12499.75  # J86
$ws.Cells.Item(86, 12).Value = 12499.75  # L86
$ws.Cells.Item(86, 14).Value = -14745.75  # N86
$ws.Cells.Item(89, 8).Value = 13858  # H89
$ws.Cells.Item(89, 10).Value = 12499.75  # J89
$ws.Cells.Item(89, 12).Value = 62498.75  # L89
$ws.Cells.Item(89, 14).Value = -73730.75  # N89
$ws.Cells.Item(115, 8).Value = 58666.332  # H115
$ws.Cells.Item(115, 9).Value = 65499.5  # I115
$ws.Cells.Item(115, 10).Value = 45000  # J115
$ws.Cells.Item(115, 11).Value = 65499.5  # K115
$ws.Cells.Item(115, 12).Value = 45000  # L115
$ws.Cells.Item(115, 13).Value = -64324.5  # M115
$ws.Cells.Item(115, 14).Value = -47350  # N115
$ws.Cells.Item(120, 8).Value = 30724.916  # H120
$ws.Cells.Item(120, 9).Value = 32496  # I120
$ws.Cells.Item(120, 10).Value = 30370.7  # J120
$ws.Cells.Item(120, 11).Value = 32496  # K120
$ws.Cells.Item(120, 12).Value = 30370.7  # L120
$ws.Cells.Item(120, 13).Value = -28867  # M120
$ws.Cells.Item(120, 14).Value = -37628.7  # N120
$ws.Cells.Item(133, 8).Value = 63607.816  # H133
$ws.Cells.Item(133, 10).Value = 66668.60000000001  # J133
$ws.Cells.Item(133, 12).Value = 66668.60000000001  # L133
$ws.Cells.Item(133, 14).Value = -71728.60000000001  # N133
$ws.Cells.Item(141, 8).Value = 102995.5  # H141
$ws.Cells.Item(141, 10).Value = 102995.5  # J141
$ws.Cells.Item(141, 12).Value = 102995.5  # L141
$ws.Cells.Item(141, 14).Value = -113355.5  # N141

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(128, 8).Value = 197014.5  # H128
$ws.Cells.Item(128, 9).Value = 197014.5  # I128
$ws.Cells.Item(128, 11).Value = 591043.5  # K128
$ws.Cells.Item(128, 13).Value = -586063.5  # M128
$ws.Cells.Item(140, 8).Value = 3327.6667  # H140
$ws.Cells.Item(140, 9).Value = 3466.6667  # I140
$ws.Cells.Item(140, 10).Value = 3188.6667  # J140
$ws.Cells.Item(140, 11).Value = 10400.0001  # K140
$ws.Cells.Item(140, 12).Value = 9566.000100000001  # L140
$ws.Cells.Item(140, 13).Value = -5220.000100000001  # M140
$ws.Cells.Item(140, 14).Value = -19926.0001  # N140

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(47, 8).Value = 29950  # H47
$ws.Cells.Item(47, 10).Value = 29950  # J47
$ws.Cells.Item(47, 12).Value = 29950  # L47
$ws.Cells.Item(47, 14).Value = -31086  # N47
$ws.Cells.Item(52, 8).Value = 34974.75  # H52
$ws.Cells.Item(52, 10).Value = 34999.668  # J52
$ws.Cells.Item(52, 12).Value = 34999.668  # L52
$ws.Cells.Item(52, 14).Value = -35517.668  # N52
$ws.Cells.Item(62, 8).Value = 72000  # H62
$ws.Cells.Item(62, 9).Value = 72000  # I62
$ws.Cells.Item(62, 11).Value = 72000  # K62
$ws.Cells.Item(62, 13).Value = -71314  # M62
$ws.Cells.Item(63, 8).Value = 40000  # H63
$ws.Cells.Item(63, 10).Value = 0  # J63
$ws.Cells.Item(63, 12).Value = 0  # L63
$ws.Cells.Item(63, 14).ClearContents()  # N63
$ws.Cells.Item(64, 8).Value = 0  # H64
$ws.Cells.Item(64, 10).Value = 0  # J64
$ws.Cells.Item(64, 12).Value = 0  # L64
$ws.Cells.Item(64, 14).ClearContents()  # N64
$ws.Cells.Item(65, 8).Value = 72000  # H65
$ws.Cells.Item(65, 9).Value = 72000  # I65
$ws.Cells.Item(65, 11).Value = 216000  # K65
$ws.Cells.Item(65, 13).Value = -212568  # M65
$ws.Cells.Item(66, 8).Value = 40000  # H66
$ws.Cells.Item(66, 10).Value = 0  # J66
$ws.Cells.Item(66, 12).Value = 0  # L66
$ws.Cells.Item(66, 14).ClearContents()  # N66
$ws.Cells.Item(67, 8).Value = 0  # H67
$ws.Cells.Item(67, 10).Value = 0  # J67
$ws.Cells.Item(67, 12).Value = 0  # L67
$ws.Cells.Item(67, 14).ClearContents()  # N67
$ws.Cells.Item(122, 8).Value = 16668909  # H122
$ws.Cells.Item(122, 9).Value = 2230.1177  # I122
$ws.Cells.Item(122, 10).Value = 38463796  # J122
$ws.Cells.Item(122, 11).Value = 6690.353099999999  # K122
$ws.Cells.Item(122, 12).Value = 115391388  # L122
$ws.Cells.Item(122, 13).Value = -4240.353099999999  # M122
$ws.Cells.Item(122, 14).Value = -115396288  # N122

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(6, 8).Value = 48240  # H6
$ws.Cells.Item(6, 10).Value = 48240  # J6
$ws.Cells.Item(6, 12).Value = 48240  # L6
$ws.Cells.Item(6, 14).Value = -48464  # N6
$ws.Cells.Item(22, 8).Value = 1727.625  # H22
$ws.Cells.Item(22, 10).Value = 1616  # J22
$ws.Cells.Item(22, 12).Value = 1616  # L22
$ws.Cells.Item(22, 14).Value = -2206  # N22
$ws.Cells.Item(27, 8).Value = 1727.625  # H27
$ws.Cells.Item(27, 10).Value = 1616  # J27
$ws.Cells.Item(27, 12).Value = 1616  # L27
$ws.Cells.Item(27, 14).Value = -1830  # N27
$ws.Cells.Item(61, 8).Value = 8316.5  # H61
$ws.Cells.Item(61, 9).Value = 5700  # I61
$ws.Cells.Item(61, 11).Value = 5700  # K61
$ws.Cells.Item(61, 13).Value = -5498  # M61
$ws.Cells.Item(80, 8).Value = 37997.5  # H80
$ws.Cells.Item(80, 9).Value = 23995  # I80
$ws.Cells.Item(80, 11).Value = 23995  # K80
$ws.Cells.Item(80, 13).Value = -22872  # M80
$ws.Cells.Item(82, 8).Value = 3370.7144  # H82
$ws.Cells.Item(82, 10).Value = 2700.75  # J82
$ws.Cells.Item(82, 12).Value = 2700.75  # L82
$ws.Cells.Item(82, 14).Value = -3422.75  # N82
$ws.Cells.Item(83, 8).Value = 37997.5  # H83
$ws.Cells.Item(83, 9).Value = 23995  # I83
$ws.Cells.Item(83, 11).Value = 71985  # K83
$ws.Cells.Item(83, 13).Value = -66369  # M83
$ws.Cells.Item(85, 8).Value = 3370.7144  # H85
$ws.Cells.Item(85, 10).Value = 2700.75  # J85
$ws.Cells.Item(85, 12).Value = 2700.75  # L85
$ws.Cells.Item(85, 14).Value = -5196.75  # N85
$ws.Cells.Item(96, 8).Value = 30000  # H96
$ws.Cells.Item(96, 10).Value = 30000  # J96
$ws.Cells.Item(96, 12).Value = 30000  # L96
$ws.Cells.Item(96, 14).Value = -35492  # N96
$ws.Cells.Item(109, 8).Value = 44974.75  # H109
$ws.Cells.Item(109, 10).Value = 34999.5  # J109
$ws.Cells.Item(109, 12).Value = 34999.5  # L109
$ws.Cells.Item(109, 14).Value = -37773.5  # N109
$ws.Cells.Item(113, 8).Value = 8316.5  # H113
$ws.Cells.Item(113, 9).Value = 5700  # I113
$ws.Cells.Item(113, 11).Value = 5700  # K113
$ws.Cells.Item(113, 13).Value = -3530  # M113
$ws.Cells.Item(136, 8).Value = 1674164.1  # H136
$ws.Cells.Item(136, 9).Value = 2503996.5  # I136
$ws.Cells.Item(136, 10).Value = 14499  # J136
$ws.Cells.Item(136, 11).Value = 7511989.5  # K136
$ws.Cells.Item(136, 12).Value = 43497  # L136
$ws.Cells.Item(136, 13).Value = -7509439.5  # M136
$ws.Cells.Item(136, 14).Value = -48597  # N136

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(70, 8).Value = 32996.5  # H70
$ws.Cells.Item(70, 10).Value = 32996.5  # J70
$ws.Cells.Item(70, 12).Value = 32996.5  # L70
$ws.Cells.Item(70, 14).Value = -33626.5  # N70
$ws.Cells.Item(73, 8).Value = 32996.5  # H73
$ws.Cells.Item(73, 10).Value = 32996.5  # J73
$ws.Cells.Item(73, 12).Value = 32996.5  # L73
$ws.Cells.Item(73, 14).Value = -35180.5  # N73
$ws.Cells.Item(132, 8).Value = 2105.0715  # H132
$ws.Cells.Item(132, 9).Value = 1920.4615  # I132
$ws.Cells.Item(132, 11).Value = 5761.3845  # K132
$ws.Cells.Item(132, 13).Value = -3231.3845  # M132
$ws.Cells.Item(133, 8).Value = 71342.8  # H133
$ws.Cells.Item(133, 10).Value = 71342.8  # J133
$ws.Cells.Item(133, 12).Value = 71342.8  # L133
$ws.Cells.Item(133, 14).Value = -81462.8  # N133

